$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '42.762.34'
    'E2' = '  -5.39%  '
    'D3' = '2.211.06'
    'E3' = '  -6.48%  '
    'E4' = '  -0.03%  '
    'D5' = '''316.08'
    'E5' = '  +1.38%  '
    'D6' = '''99.38'
    'E6' = '  -8.54%  '
    'D7' = '''0.589'
    'E7' = '  -6.54%  '
    'E8' = '  -0.13%  '
    'D9' = '''0.561'
    'E9' = '  -8.09%  '
    'D10' = '''36.83'
    'E10' = '  -9.74%  '
    'D11' = '''54.06'
    'E11' = '  -2.91%  '
    'E12' = '  -9.58%  '
    'D13' = '''7.72'
    'E13' = '  -8.63%  '
    'D14' = '''0.106'
    'E14' = '  -2.75%  '
    'E15' = '  -11.90%  '
    'D16' = '2.544.67'
    'E16' = '  -6.59%  '
    'D17' = '''14.17'
    'D18' = '2.212.10'
    'E18' = '  -6.32%  '
    'D19' = '42.731.24'
    'E19' = '  -5.31%  '
    'D20' = '''14.69'
    'E20' = '  +2.36%  '
    'E21' = '  -9.37%  '
    'D22' = '''6.41'
    'E22' = '  -10.87%  '
    'D23' = '''65.23'
    'E23' = '  -10.74%  '
    'E24' = '  -9.85%  '
    'D25' = '''236.12'
    'E25' = '  -8.76%  '
    'D26' = '''2.13'
    'E26' = '  -7.57%  '
    'E27' = '  +0.26%  '
    'D28' = '''10.01'
    'E29' = '  -4.50%  '
    'D30' = '''6.32'
    'E30' = '  -12.28%  '
    'D31' = '''0.0895'
    'E31' = '  -7.91%  '
    'D32' = '''20.53'
    'E32' = '  -8.05%  '
    'D33' = '''34.25'
    'E33' = '  -7.67%  '
    'D34' = '''155.81'
    'E34' = '  -7.21%  '
    'E35' = '  -7.03%  '
    'E36' = '  +9.88%  '
    'E37' = '  +9.05%  '
    'E38' = '  -6.78%  '
    'B39' = 'Kaspa'
    'C39' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D39' = '''0.107'
    'E39' = '  -8.40%  '
    'B40' = 'RenderToken'
    'C40' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D40' = '''4.42'
    'E40' = '  -5.28%  '
    'B41' = 'NEARProtocol'
    'C41' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D41' = '''3.85'
    'E41' = '  -2.83%  '
    'D42' = '''0.0325'
    'E42' = '  -7.88%  '
    'D43' = '1.893.46'
    'E43' = '  +1.33%  '
    'E44' = '  +0.16%  '
    'D45' = '''12.53'
    'E45' = '  -2.39%  '
    'D46' = '''87.94'
    'E46' = '  -11.46%  '
    'D47' = '''0.208'
    'E47' = '  -9.06%  '
    'D48' = '''61.20'
    'E48' = '  -11.74%  '
    'E49' = '  -5.60%  '
    'D50' = '''76.40'
    'E50' = '  -8.46%  '
    'D51' = '''102.27'
    'E51' = '  -7.14%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
